# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.149.59"
$ws.Range("E2").Value = "  -0.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.39"
$ws.Range("E3").Value = "  -0.96%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  -1.11%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.42"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("E11").Value = "  +0.15%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.643.37"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("E13").Value = "  -0.42%  "

$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "64.84"
$ws.Range("E15").Value = "  -3.74%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.141.29"
$ws.Range("E16").Value = "  -0.01%  "

$ws.Range("E17").Value = "  +0.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.29"
$ws.Range("E18").Value = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.00"
$ws.Range("E19").Value = "  -0.18%  "

$ws.Range("E20").Value = "  +1.63%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.44"
$ws.Range("E22").Value = "  -5.23%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.08"
$ws.Range("E23").Value = "  -1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "147.83"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("E25").Value = "  -0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.32"
$ws.Range("E26").Value = "  -2.94%  "

$ws.Range("E27").Value = "  -0.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.63"
$ws.Range("E28").Value = "  -0.86%  "

$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("E30").Value = "  -1.09%  "

$ws.Range("E31").Value = "  -0.59%  "

$ws.Range("E32").Value = "  -1.00%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.349.37"
$ws.Range("E33").Value = "  +6.51%  "

$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0178"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("E38").Value = "  -0.18%  "

$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("E40").Value = "  -0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.54"
$ws.Range("E42").Value = "  +5.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.25"
$ws.Range("E43").Value = "  -0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.765.25"
$ws.Range("E44").Value = "  -1.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.63"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.853"
$ws.Range("E47").Value = "  +27.97%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1000"
$ws.Range("E48").Value = "  +2.42%  "

$ws.Range("E49").Value = "  +0.15%  "

$ws.Range("E50").Value = "  -0.67%  "

$ws.Range("E51").Value = "  -0.13%  "

